$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 29

$ws.Range("D9").Value = 19

$ws.Range("A10").Value = "asd"
$ws.Range("B10").Value = "asd"
$ws.Range("C10").Value = "F9B351E9"
$ws.Range("D10").Value = 18
$ws.Range("E10").Formula = '="101"'
$ws.Range("F10").Formula = '="TRUE"'
$ws.Range("E10:F10").Copy()
$ws.Range("E10:F10").PasteSpecial(-4163)
